$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# SYSC 4504 (sheet1) - Lab 6 deliverable grade entered
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("SYSC 4504")
$ws1.Range("D10").Value = 1

# ---------------------------------------------------------------------------
# SYSC 4502 (sheet2) - Assignment 3 grade updated + self-predicted final grade
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("SYSC 4502")
$ws2.Range("D16").Value = 1
$ws2.Range("D29").Value = 0.83

# ---------------------------------------------------------------------------
# SYSC 4415 (sheet3) - weekly reflection grades, grades, self-predicted grades
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("SYSC 4415")
$ws3.Range("E5").Value = 0.06
$ws3.Range("E6").Value = 0.06
$ws3.Range("D7").Value = 1
$ws3.Range("E7").Value = 0.08
$ws3.Range("D23").Value = 0.93
$ws3.Range("D29").Value = 0.92
$ws3.Range("D33").Value = 0.89

# ---------------------------------------------------------------------------
# ECOR 4995 (sheet4) - deliverable grade, total formula fix, final exam tweak,
# self-predicted grade table added
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("ECOR 4995")
$ws4.Range("D8").Value = 1
$ws4.Range("F9").Formula = "=SUM(F5:F8)"
$ws4.Range("D23").Value = 1
$ws4.Range("F30").Formula = "=SUM(F26*E26,D27*E27)"

$ws4.Range("C32").Value = "Get:"
$ws4.Range("D32").Value = 0.88
$ws4.Range("E32").Value = "A+"
$ws4.Range("D33").Value = 0.79
$ws4.Range("E33").Value = "A"
$ws4.Range("D34").Value = 0.71
$ws4.Range("E34").Value = "A-"

# ---------------------------------------------------------------------------
# View-state: SYSC 4415 becomes the active/selected sheet
# ---------------------------------------------------------------------------
$ws3.Activate()
$ws3.Range("A1:XFD1048576").Select()
